# Auto-generated Excel COM-interop script
# Applies literal cell value updates per the target diff (Mateus_Profits workbook).
# All target cells are plain numeric, cached (non-formula) values pulled from an external
# market-data source; some cells are cleared entirely (no HQ/NQ price -> profit cell removed)
# and one cell is newly populated (HQ price now available).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4442.1
$ws.Range("J19").Value = 4227.7144
$ws.Range("L19").Value = 4227.7144
$ws.Range("N19").Value = -4577.7144
$ws.Range("H58").Value = 634.8421
$ws.Range("J58").Value = 804.44446
$ws.Range("L58").Value = 2413.33338
$ws.Range("N58").Value = -2713.33338
$ws.Range("H62").Value = 5102.2
$ws.Range("I62").Value = 5429.533
$ws.Range("J62").Value = 4120.2
$ws.Range("K62").Value = 5429.533
$ws.Range("L62").Value = 4120.2
$ws.Range("M62").Value = -4805.533
$ws.Range("N62").Value = -5368.2
$ws.Range("H65").Value = 5102.2
$ws.Range("I65").Value = 5429.533
$ws.Range("J65").Value = 4120.2
$ws.Range("K65").Value = 27147.665
$ws.Range("L65").Value = 20601
$ws.Range("M65").Value = -24027.665
$ws.Range("N65").Value = -26841
$ws.Range("H70").Value = 3075.5
$ws.Range("I70").Value = 1706
$ws.Range("J70").Value = 4445
$ws.Range("K70").Value = 5118
$ws.Range("L70").Value = 13335
$ws.Range("M70").Value = -4848
$ws.Range("N70").Value = -13875
$ws.Range("H73").Value = 3075.5
$ws.Range("I73").Value = 1706
$ws.Range("J73").Value = 4445
$ws.Range("K73").Value = 5118
$ws.Range("L73").Value = 13335
$ws.Range("M73").Value = -4182
$ws.Range("N73").Value = -15207
$ws.Range("H74").Value = 10696.923
$ws.Range("I74").Value = 8443.666999999999
$ws.Range("K74").Value = 8443.666999999999
$ws.Range("M74").Value = -7507.666999999999
$ws.Range("H77").Value = 10696.923
$ws.Range("I77").Value = 8443.666999999999
$ws.Range("K77").Value = 42218.335
$ws.Range("M77").Value = -37538.335
$ws.Range("H86").Value = 2999
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 2999
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H132").Value = 1050.8
$ws.Range("I132").Value = 1050.8
$ws.Range("K132").Value = 3152.4
$ws.Range("M132").Value = -622.3999999999996
$ws.Range("H137").Value = 4872.2144
$ws.Range("I137").Value = 6173.5625
$ws.Range("K137").Value = 18520.6875
$ws.Range("M137").Value = -15970.6875
$ws.Range("H138").Value = 8136.6787
$ws.Range("I138").Value = 11656.417
$ws.Range("J138").Value = 5496.875
$ws.Range("K138").Value = 34969.251
$ws.Range("L138").Value = 16490.625
$ws.Range("M138").Value = -29829.251
$ws.Range("N138").Value = -26770.625
$ws.Range("H141").Value = 1621.138
$ws.Range("I141").Value = 1621.138
$ws.Range("K141").Value = 4863.414
$ws.Range("M141").Value = 316.5860000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13639.821
$ws.Range("I32").Value = 9876.08
$ws.Range("K32").Value = 9876.08
$ws.Range("M32").Value = -9589.08
$ws.Range("H61").Value = 4623
$ws.Range("I61").Value = 4623
$ws.Range("K61").Value = 4623
$ws.Range("M61").Value = -4411
$ws.Range("H136").Value = 4623
$ws.Range("I136").Value = 4623
$ws.Range("K136").Value = 13869
$ws.Range("M136").Value = -11319

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 18510
$ws.Range("I97").Value = 16769.75
$ws.Range("K97").Value = 16769.75
$ws.Range("M97").Value = -15778.75
$ws.Range("H105").Value = 2402.7144
$ws.Range("I105").Value = 2402.7144
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2402.7144
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -655.7143999999998
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 44951.715
$ws.Range("J106").Value = 44951.715
$ws.Range("L106").Value = 44951.715
$ws.Range("N106").Value = -47475.715
$ws.Range("H134").Value = 6311
$ws.Range("I134").Value = 3671.5557
$ws.Range("K134").Value = 11014.6671
$ws.Range("M134").Value = -8479.667099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 56817.145
$ws.Range("J18").Value = 56817.145
$ws.Range("L18").Value = 56817.145
$ws.Range("N18").Value = -57277.145
$ws.Range("H31").Value = 5303.8
$ws.Range("I31").Value = 3198.9092
$ws.Range("J31").Value = 7876.4443
$ws.Range("K31").Value = 3198.9092
$ws.Range("L31").Value = 7876.4443
$ws.Range("M31").Value = -2903.9092
$ws.Range("N31").Value = -8466.444299999999
$ws.Range("H34").Value = 5303.8
$ws.Range("I34").Value = 3198.9092
$ws.Range("J34").Value = 7876.4443
$ws.Range("K34").Value = 3198.9092
$ws.Range("L34").Value = 7876.4443
$ws.Range("M34").Value = -2996.9092
$ws.Range("N34").Value = -8280.444299999999
$ws.Range("H41").Value = 20449.5
$ws.Range("I41").Value = 20449.5
$ws.Range("K41").Value = 20449.5
$ws.Range("M41").Value = -20021.5
$ws.Range("H110").Value = 79573
$ws.Range("J110").Value = 79573
$ws.Range("L110").Value = 79573
$ws.Range("N110").Value = -87753
$ws.Range("H117").Value = 28548.666
$ws.Range("J117").Value = 28548.666
$ws.Range("L117").Value = 28548.666
$ws.Range("N117").Value = -37726.666
$ws.Range("H132").Value = 268082.47
$ws.Range("I132").Value = 338103.62
$ws.Range("K132").Value = 1014310.86
$ws.Range("M132").Value = -1011780.86
$ws.Range("H141").Value = 512262.9
$ws.Range("J141").Value = 608988
$ws.Range("L141").Value = 608988
$ws.Range("N141").Value = -619348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 600
$ws.Range("I9").Value = 600
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 1800
$ws.Range("L9").Value = 1800
$ws.Range("M9").Value = -1576
$ws.Range("N9").Value = -2248
$ws.Range("H12").Value = 687.7143
$ws.Range("J12").Value = 1055
$ws.Range("L12").Value = 3165
$ws.Range("N12").Value = -3511
$ws.Range("H131").Value = 5895.6665
$ws.Range("J131").Value = 6158.909
$ws.Range("L131").Value = 18476.727
$ws.Range("N131").Value = -28556.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2501000.5
$ws.Range("J3").Value = 5000001
$ws.Range("L3").Value = 5000001
$ws.Range("N3").Value = -5000233
$ws.Range("H14").Value = 11279.333
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H18").Value = 6264999.5
$ws.Range("J18").Value = 19999.334
$ws.Range("L18").Value = 19999.334
$ws.Range("N18").Value = -20585.334
$ws.Range("H24").Value = 13807.25
$ws.Range("I24").Value = 416
$ws.Range("J24").Value = 15720.286
$ws.Range("K24").Value = 416
$ws.Range("L24").Value = 15720.286
$ws.Range("M24").Value = -243
$ws.Range("N24").Value = -16066.286
$ws.Range("H80").Value = 3450
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 3450
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 184205.27
$ws.Range("I61").Value = 287651.44
$ws.Range("K61").Value = 287651.44
$ws.Range("M61").Value = -287449.44
$ws.Range("H113").Value = 184205.27
$ws.Range("I113").Value = 287651.44
$ws.Range("K113").Value = 287651.44
$ws.Range("M113").Value = -285481.44
$ws.Range("H132").Value = 332993.1
$ws.Range("I132").Value = 352192.7
$ws.Range("K132").Value = 1056578.1
$ws.Range("M132").Value = -1054048.1
$ws.Range("H136").Value = 733333300
$ws.Range("I136").Value = 600000000
$ws.Range("K136").Value = 1800000000
$ws.Range("M136").Value = -1799997450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 55866.156
$ws.Range("I107").Value = 2652.9092
$ws.Range("J107").Value = 129034.375
$ws.Range("K107").Value = 7958.7276
$ws.Range("L107").Value = 387103.125
$ws.Range("M107").Value = -6038.7276
$ws.Range("N107").Value = -390943.125
$ws.Range("H113").Value = 816.4737
$ws.Range("J113").Value = 399.33334
$ws.Range("L113").Value = 1198.00002
$ws.Range("N113").Value = -5538.000019999999
$ws.Range("H140").Value = 125000
$ws.Range("J140").Value = 125000
$ws.Range("L140").Value = 125000
$ws.Range("N140").Value = -135360
$ws.Range("H141").Value = 80714.664
$ws.Range("J141").Value = 80714.664
$ws.Range("L141").Value = 80714.664
$ws.Range("N141").Value = -91074.664
